# Append the new daily allocation row (11/14/2025) to the bottom of the
# sheet, mirroring the existing rows (Date text in col A, fractional
# allocations in cols B/C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 74

# Force column A to be treated as plain text so the date-like string isn't
# auto-converted into a date serial number, then drop the format override
# so the new cell ends up unstyled (matching the rest of the data rows).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "11/14/2025"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = 0.2019497281731677
$ws.Cells.Item($row, 3).Value = 0.7980502718268323
